$wb = $excel.ActiveWorkbook

# Sheet1: B column counter-style values change (A column unchanged)
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("B2").Value = 1
$ws1.Range("B3").Value = 2
$ws1.Range("B4").Value = 1
$ws1.Range("B5").Value = 1
$ws1.Range("B6").Value = 1

# Sheet2: A4/A6 flags flip, B column values change
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("B3").Value = 2
$ws2.Range("A4").Value = 0
$ws2.Range("B4").Value = 2
$ws2.Range("B5").Value = 2
$ws2.Range("A6").Value = 1
$ws2.Range("B6").Value = 2

# Sheet5: A2/A3 flags flip, B column values change (B5 gains a value)
$ws5 = $wb.Worksheets.Item("Sheet5")
$ws5.Range("A2").Value = 1
$ws5.Range("B2").Value = 1
$ws5.Range("A3").Value = 0
$ws5.Range("B3").Value = 2
$ws5.Range("B4").Value = 3
$ws5.Range("B5").Value = 4
$ws5.Range("B6").Value = 5
